# Sample model workbook: add "Reaction List" / "Metabolite List" data
# (mirrors the authored diff: rename Sheet1, add a metabolite sheet,
# populate both with header + sample rows, light header formatting,
# and leave the Metabolite List tab active/selected).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Sheets: rename Sheet1 -> "Reaction List", add "Metabolite List"
# ---------------------------------------------------------------
$wsReactions = $wb.Worksheets.Item(1)
$wsReactions.Name = "Reaction List"

$wsMetabolites = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsReactions)
$wsMetabolites.Name = "Metabolite List"

# ---------------------------------------------------------------
# 2. Reaction List sheet
# ---------------------------------------------------------------
$reactionHeaders = @("Abbreviation", "Description", "Reaction", "GPR", "Genes", "Proteins", "Subsystem", "Reversible", "Lower bound", "Upper bound", "Objective", "Confidence", "EC Number", "Notes", "References")
for ($col = 0; $col -lt $reactionHeaders.Length; $col++) {
    $wsReactions.Cells.Item(1, $col + 1).Value = $reactionHeaders[$col]
}
$wsReactions.Rows.Item(1).RowHeight = 15

# Lower bound / Upper bound headers pick up a (no-op) fill format flag
$wsReactions.Range("I1:J1").Interior.ColorIndex = -4142

$abbrevs = @("R1", "R2", "R3", "R4")
$reactions = @(" -> A", " -> B", "A + B -> C", "C ->")
$reversible = @($false, $true, $false, $false)
$lower = @(0, -10, 0, 0)
$upper = @(10, 10, 10, 10)
$objective = @(0, 0, 0, 1)

for ($i = 0; $i -lt 4; $i++) {
    $r = $i + 2
    $wsReactions.Cells.Item($r, 1).Value = $abbrevs[$i]
    $wsReactions.Cells.Item($r, 3).Value = $reactions[$i]
    $wsReactions.Cells.Item($r, 8).Value = $reversible[$i]
    $wsReactions.Cells.Item($r, 9).Value = $lower[$i]
    $wsReactions.Cells.Item($r, 10).Value = $upper[$i]
    $wsReactions.Cells.Item($r, 11).Value = $objective[$i]
}

# ---------------------------------------------------------------
# 3. Metabolite List sheet
# ---------------------------------------------------------------
$metaboliteHeaders = @("Abbreviation", "Description", "Neutral formula", "Charged formula", "Charge", "Compartment", "KEGG ID", "PubChem ID", "ChEBI ID", "InChI string", "SMILES", "HMDB ID")
for ($col = 0; $col -lt $metaboliteHeaders.Length; $col++) {
    $wsMetabolites.Cells.Item(1, $col + 1).Value = $metaboliteHeaders[$col]
}

$wsMetabolites.Cells.Item(2, 1).Value = "A"
$wsMetabolites.Cells.Item(3, 1).Value = "B"
$wsMetabolites.Cells.Item(4, 1).Value = "C"

# ---------------------------------------------------------------
# 4. Selection / active tab: Reaction List keeps K2 selected,
#    Metabolite List (A3 selected) ends up the active sheet.
# ---------------------------------------------------------------
$wsReactions.Range("K2").Select() | Out-Null
$wsMetabolites.Range("A3").Select() | Out-Null
